# Update countries & provincias Spain
# Refreshes the COVID "Pais" data snapshot:
#   - timestamp string bumped from 03:25 to 04:42
#   - Bolivia re-sorted above Barein/Israel (new case counts)
#   - Australia / Paraguay case counts refreshed
#   - small-territory rows (Groenlandia/Islas Malvinas,
#     Montserrat/Seychelles, Islas Virgenes Britanicas/Papua Nueva
#     Guinea) swap order per the refreshed ranking

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 17 de Junio de 2020 a las 04:42"

# --- Rows 49-51: Bolivia moves above Barein/Israel with fresh data ----
$ws.Range("A49").Value = "Bolivia"
$ws.Range("B49").Value = 19883
$ws.Range("C49").Value = 810
$ws.Range("D49").Value = 3752
$ws.Range("E49").Value = 15472
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 27
$ws.Range("H49").Value = 659

$ws.Range("A50").Value = "Barein"
$ws.Range("B50").Value = 19553
$ws.Range("C50").Value = 0
$ws.Range("D50").Value = 13866
$ws.Range("E50").Value = 5640
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 47

$ws.Range("A51").Value = "Israel"
$ws.Range("B51").Value = 19495
$ws.Range("C51").Value = 0
$ws.Range("D51").Value = 15449
$ws.Range("E51").Value = 3744
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 302

# --- Row 72: Australia data refresh ------------------------------------
$ws.Range("B72").Value = 7370
$ws.Range("C72").Value = 23
$ws.Range("D72").Value = 6859
$ws.Range("E72").Value = 409

# --- Row 121: Paraguay data refresh ------------------------------------
$ws.Range("E121").Value = 591
$ws.Range("G121").Value = 1
$ws.Range("H121").Value = 13

# --- Rows 206-207: Islas Malvinas / Groenlandia swap order -------------
$ws.Range("A206").Value = "Groenlandia"
$ws.Range("A207").Value = "Islas Malvinas"

# --- Rows 210-211: Montserrat / Seychelles swap order + data -----------
$ws.Range("A210").Value = "Seychelles"
$ws.Range("D210").Value = 11
$ws.Range("H210").Value = 0

$ws.Range("A211").Value = "Montserrat"
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1

# --- Rows 213-214: Islas Virgenes Britanicas / Papua Nueva Guinea swap -
$ws.Range("A213").Value = "Papua Nueva Guinea"
$ws.Range("D213").Value = 8
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Islas Virgenes Britanicas"
$ws.Range("D214").Value = 7
$ws.Range("H214").Value = 1
